# Update Vignola report through 20/09/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newData = @(
    @(375, 44449, 2,  22, 86.10904536381071),
    @(376, 44450, 3,  20, 78.28095033073701),
    @(377, 44451, 10, 26, 101.7652354299581),
    @(378, 44452, 4,  30, 117.4214254961055),
    @(379, 44453, 7,  35, 136.9916630787898),
    @(380, 44454, 0,  34, 133.0776155622529),
    @(381, 44455, 9,  35, 136.9916630787898),
    @(382, 44456, 1,  34, 133.0776155622529),
    @(383, 44457, 6,  37, 144.8197581118635),
    @(384, 44458, 9,  36, 140.9057105953266),
    @(385, 44459, 2,  34, 133.0776155622529)
)

# Copy the date-column format from the last existing row (A374) so new date
# cells (column A) keep the same style (border, bold, centered, date format).
$ws.Range("A374").Copy()

foreach ($row in $newData) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]

    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 1)).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

Write-Host "Updated rows 375-385 on sheet $($ws.Name)"
